$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet2 ("建物"): row4 share_portion "全部■" -> "全部"
# ---------------------------------------------------------------------------
$wsBuilding = $wb.Worksheets.Item(2)
$wsBuilding.Cells.Item(4,4).Value = "全部"

# ---------------------------------------------------------------------------
# 2) Sheet7 ("保險"): strip the trailing "■" from the insurance period text
# ---------------------------------------------------------------------------
$wsInsurance = $wb.Worksheets.Item(7)
$wsInsurance.Cells.Item(1,5).Value = "缴費期滿保險期間：終身"
$wsInsurance.Cells.Item(2,5).Value = "缴費期滿保險期間：終身"

# ---------------------------------------------------------------------------
# 3) Sheet6 ("具有相當價值之財產"): turn row 1 into a real header row and add
#    the property_category/category/date/legislator_name/legislator_id/
#    source_file/index columns (F:L) to every data row.
# ---------------------------------------------------------------------------
$wsAsset = $wb.Worksheets.Item(6)

# Copy formatting (bold header style) from B1 across to F1:L1 before
# overwriting the values, then fix up B1:E1 to hold the column headers.
$wsAsset.Range("B1").Copy($wsAsset.Range("F1:L1"))

$wsAsset.Cells.Item(1,2).Value = "name"
$wsAsset.Cells.Item(1,3).Value = "quantity"
$wsAsset.Cells.Item(1,4).Value = "owner"
$wsAsset.Cells.Item(1,5).Value = "total"
$wsAsset.Cells.Item(1,6).Value = "property_category"
$wsAsset.Cells.Item(1,7).Value = "category"
$wsAsset.Cells.Item(1,8).Value = "date"
$wsAsset.Cells.Item(1,9).Value = "legislator_name"
$wsAsset.Cells.Item(1,10).Value = "legislator_id"
$wsAsset.Cells.Item(1,11).Value = "source_file"
$wsAsset.Cells.Item(1,12).Value = "index"

# Copy formatting (plain data style) from B2 across to F2:L6, then fill in
# the values for every data row (2-6).
$wsAsset.Range("B2").Copy($wsAsset.Range("F2:L6"))

$rows = 2..6
$indexValues = @{2=83;3=84;4=85;5=86;6=87}
foreach ($r in $rows) {
    $wsAsset.Cells.Item($r,6).Value = "otherbonds"
    $wsAsset.Cells.Item($r,7).Value = "normal"
    $wsAsset.Cells.Item($r,8).Value = "2011-11-21"
    $wsAsset.Cells.Item($r,9).Value = "孫大千"
    $wsAsset.Cells.Item($r,10).Value = 919
    $wsAsset.Cells.Item($r,11).Value = "tmpc6841"
    $wsAsset.Cells.Item($r,12).Value = $indexValues[$r]
}

Write-Host "edit complete"
